# Re-orders the sample-request rows in the CSDCO dir-walk report.
# The "A" (site) column never changes, but B/C/D travel together as a
# logical record that gets shuffled to new row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/~`$hnson_SRF_part1of2.doc"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/GLAD7JohnsonTEX86.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF2_GLAD7_Scholz_20070112.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/images/GLAD7-MAL05-1B-45E-2 copy.bmp"; C = "Images"; D = "valid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/images/GLAD7-MAL05-1B-43E-2 copy.bmp"; C = "Images"; D = "valid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/Johnson_SRF_part2of2.doc"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "GLAD7-metadata.xls"; C = "metadata"; D = "valid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/images/GLAD7-MAL05-1B-44E-4 copy.bmp"; C = "Images"; D = "valid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/GLAD7_Petrick_C14_20090302/LacCore_GLAD7_SRF_part1of2.doc"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/images/GLAD7-MAL05-1B-42E-3 copy.bmp"; C = "Images"; D = "valid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/LacCore_SRF_part2of2_Abbott_20090917.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/MAL05-TCJ-sampls.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/GLAD7_Stone_diatoms_20061026/GLAD7_Stone_diatoms_20061026.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/MALsamples4_15to4_17.xlsx"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/LacCore_SRF_part2of2_GLAD7_Johnson_20081201.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/GLAD7_Beuning_pollen_20061103/LacCore_SRF_part2of2 Nov 2006.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/MW LacCore_SRF_part2of2.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/GLAD7_Petrick_C14_20090302/LacCore_GLAD7_SRF_part2of2 ETB.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/Johnson_SRF_part1of2.doc"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/images/GLAD7-MAL05-1B-42E-2 copy.bmp"; C = "Images"; D = "valid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/GLAD7 sampling_20070109.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/LacCore_SRF_part2of2complete_Ivory.xls"; C = ""; D = "notvalid" },
    @{ A = "GLAD7 (Malawi)"; B = "/sample requests/SRF_GLAD7_Johnson_20081201/ptf-Johnson_GLAD7_20081201.doc"; C = ""; D = "notvalid" }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4))
    if ($row.D -eq "notvalid") {
        $rowRange.Interior.Color = 8036607
    } else {
        $rowRange.ClearFormats()
    }
}
